$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("As static variable") gains marks for the newly-tested Pointer and
# Array contexts (same plain "Y" marker used by the other columns), plus a
# distinct "x" marker (highlighted in yellow) for Struct and Union, which are
# only partially supported as static variables.
$ws.Range("H3").Value = "Y"
$ws.Range("I3").Value = "Y"

$ws.Range("J3").Value = "x"
$ws.Range("J3").Interior.Color = 65535

$ws.Range("S3").Value = "x"
$ws.Range("S3").Interior.Color = 65535

# Move the active cell selection as recorded in the workbook view.
$ws.Range("I4").Select() | Out-Null
